$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.854221333333333
$ws.Range("H2").Value = 5.562664
$ws.Range("I2").Value = 0.03416002559055492
$ws.Range("J2").Value = 0.03416002559055492
$ws.Range("M2").Value = 0.1112926666666667
$ws.Range("N2").Value = 0.333878
$ws.Range("O2").Value = 0.01397697460904174
$ws.Range("P2").Value = 0.01397697460904174
$ws.Range("Q2").Value = 0.2063612367768889
$ws.Range("R2").Value = 1.857251130992
$ws.Range("S2").Value = 0.000477453810323402
$ws.Range("T2").Value = 0.000477453810323402
$ws.Range("G3").Value = 1.854221333333333
$ws.Range("H3").Value = 5.562664
$ws.Range("I3").Value = 0.03416002559055492
$ws.Range("J3").Value = 0.03416002559055492
$ws.Range("O3").Value = 0.4165551449121381
$ws.Range("P3").Value = 0.4165551449121381
$ws.Range("Q3").Value = 6.150174647540444
$ws.Range("R3").Value = 55.351571827864
$ws.Range("S3").Value = 0.01422953441007595
$ws.Range("T3").Value = 0.01422953441007595
$ws.Range("G4").Value = 1.854221333333333
$ws.Range("H4").Value = 5.562664
$ws.Range("I4").Value = 0.03416002559055492
$ws.Range("J4").Value = 0.03416002559055492
$ws.Range("O4").Value = 0.5694678804788202
$ws.Range("P4").Value = 0.5694678804788201
$ws.Range("Q4").Value = 8.407834986285334
$ws.Range("R4").Value = 75.67051487656799
$ws.Range("S4").Value = 0.01945303737015557
$ws.Range("T4").Value = 0.01945303737015557
$ws.Range("I5").Value = 0.8311547934421808
$ws.Range("J5").Value = 0.8311547934421808
$ws.Range("M5").Value = 0.1112926666666667
$ws.Range("N5").Value = 0.333878
$ws.Range("O5").Value = 0.01397697460904174
$ws.Range("P5").Value = 0.01397697460904174
$ws.Range("Q5").Value = 5.021018812561779
$ws.Range("R5").Value = 45.18916931305601
$ws.Range("S5").Value = 0.01161702944412469
$ws.Range("T5").Value = 0.01161702944412469
$ws.Range("I6").Value = 0.8311547934421808
$ws.Range("J6").Value = 0.8311547934421808
$ws.Range("O6").Value = 0.4165551449121381
$ws.Range("P6").Value = 0.4165551449121381
$ws.Range("S6").Value = 0.3462218054267259
$ws.Range("T6").Value = 0.3462218054267259
$ws.Range("I7").Value = 0.8311547934421808
$ws.Range("J7").Value = 0.8311547934421808
$ws.Range("O7").Value = 0.5694678804788202
$ws.Range("P7").Value = 0.5694678804788201
$ws.Range("R7").Value = 1841.155270659024
$ws.Range("S7").Value = 0.4733159585713304
$ws.Range("T7").Value = 0.4733159585713302
$ws.Range("I8").Value = 0.1346851809672642
$ws.Range("J8").Value = 0.1346851809672642
$ws.Range("M8").Value = 0.1112926666666667
$ws.Range("N8").Value = 0.333878
$ws.Range("O8").Value = 0.01397697460904174
$ws.Range("P8").Value = 0.01397697460904174
$ws.Range("Q8").Value = 0.8136352370768888
$ws.Range("R8").Value = 7.322717133692
$ws.Range("S8").Value = 0.001882491354593643
$ws.Range("T8").Value = 0.001882491354593643
$ws.Range("I9").Value = 0.1346851809672642
$ws.Range("J9").Value = 0.1346851809672642
$ws.Range("O9").Value = 0.4165551449121381
$ws.Range("P9").Value = 0.4165551449121381
$ws.Range("S9").Value = 0.0561038050753363
$ws.Range("T9").Value = 0.0561038050753363
$ws.Range("I10").Value = 0.1346851809672642
$ws.Range("J10").Value = 0.1346851809672642
$ws.Range("O10").Value = 0.5694678804788202
$ws.Range("P10").Value = 0.5694678804788201
$ws.Range("S10").Value = 0.0766988845373343
$ws.Range("T10").Value = 0.07669888453733428
